# [Kadastro App] Yeni kayit eklendi: 2933
#
# Adds the new Kadastro record (Kayit No 2933, Erdemli / CAP, Sevil Saracer)
# as a new row right after the last existing row, both in the master
# "Kayitlar" list and in the unit-specific "Erdemli" sheet.
#
# The leading "'" on the number/date-looking fields forces Excel to store
# them as literal text (matching how every other row in these sheets already
# stores "Kayit No" / "Tarih" / "Parsel Sayisi" as text, not as numbers or
# dates) instead of letting Excel auto-convert them to numeric/date values.
# ClearFormats() afterwards drops the transient "quote prefix" number format
# that setting a text-look-alike value picks up, so the new cells end up
# with the same plain/General styling as the existing rows.

$wb = $excel.ActiveWorkbook

$kayitNo      = "'2933"
$tarih        = "'2025-09-08"
$birim        = "Erdemli"
$parselSayisi = "'1"
$is           = "ÇAP"
$personel     = "SEVİL SARAÇER (Tekniker)"

function Add-KayitRow {
    param($ws)

    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $kayitNo
    $ws.Cells.Item($newRow, 2).Value = $tarih
    $ws.Cells.Item($newRow, 3).Value = $birim
    $ws.Cells.Item($newRow, 4).Value = $parselSayisi
    $ws.Cells.Item($newRow, 5).Value = $is
    $ws.Cells.Item($newRow, 6).Value = $personel

    $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 6)).ClearFormats()
}

# Master list of all records
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Add-KayitRow $wsKayitlar

# Per-unit sheet for "Erdemli"
$wsErdemli = $wb.Worksheets.Item("Erdemli")
Add-KayitRow $wsErdemli
